$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values that are unambiguous as text (contain % signs, letters, URLs,
# multiple dots, leading/trailing spaces, etc.) - plain assignment is safe.
$safeUpdates = @{
    'D2' = '62.038.63'
    'D3' = '2.420.92'
    'E3' = '  -0.14%  '
    'E4' = '  +0.01%  '
    'E5' = '  -0.22%  '
    'E6' = '  -0.62%  '
    'E7' = '  -0.01%  '
    'E8' = '  -0.60%  '
    'D9' = '2.420.00'
    'E9' = '  -0.10%  '
    'E10' = '  -0.47%  '
    'E11' = '  +0.14%  '
    'E12' = '  -3.36%  '
    'E14' = '  +0.56%  '
    'E15' = '  -2.28%  '
    'D16' = '2.859.41'
    'E16' = '  -0.09%  '
    'D17' = '61.964.53'
    'E17' = '  -0.23%  '
    'D18' = '2.420.09'
    'E18' = '  -0.14%  '
    'E19' = '  -1.01%  '
    'E20' = '  -0.21%  '
    'E21' = '  -1.72%  '
    'E22' = '  +0.97%  '
    'E23' = '  -0.08%  '
    'E24' = '  +2.66%  '
    'E25' = '  +1.54%  '
    'E26' = '  -2.33%  '
    'E27' = '  -4.62%  '
    'D28' = '2.541.35'
    'E28' = '  +0.42%  '
    'E29' = '  -0.20%  '
    'D30' = '0.0₃0933'
    'E30' = '  -1.31%  '
    'E31' = '  -0.86%  '
    'E32' = '  -5.13%  '
    'E33' = '  -1.98%  '
    'E34' = '  -0.98%  '
    'E35' = '  -3.28%  '
    'E36' = '  -0.02%  '
    'B37' = 'PolygonEcosystemToken'
    'C37' = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
    'E37' = '  -1.02%  '
    'B38' = 'NEARProtocol'
    'C38' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'E38' = '  -1.57%  '
    'E39' = '  -4.73%  '
    'E40' = '  -0.77%  '
    'E41' = '  -0.08%  '
    'E42' = '  -1.72%  '
    'E43' = '  +0.49%  '
    'E44' = '  -2.74%  '
    'E45' = '  -1.98%  '
    'E46' = '  -0.91%  '
    'E47' = '  -1.49%  '
    'E48' = '  -2.44%  '
    'E49' = '  +0.00%  '
    'E51' = '  -0.46%  '
}

# Cell values that look like plain decimal numbers (single dot, digits only).
# Excel's COM layer auto-converts these to numeric values on assignment, so
# force the cell to Text format first to preserve them as literal strings,
# exactly like the original "Price" column entries.
$textUpdates = @{
    'D5' = '562.51'
    'D6' = '143.66'
    'D8' = '0.529'
    'D13' = '0.349'
    'D14' = '26.13'
    'D15' = '0.0000173'
    'D19' = '11.24'
    'D20' = '323.06'
    'D21' = '4.13'
    'D24' = '67.33'
    'D26' = '8.71'
    'D27' = '557.03'
    'D31' = '8.21'
    'D37' = '0.380'
    'D38' = '4.73'
    'D40' = '152.41'
    'D41' = '18.67'
    'D44' = '2.27'
    'D45' = '147.44'
    'D47' = '0.0530'
    'D48' = '19.89'
    'D50' = '0.0921'
}

foreach ($addr in $safeUpdates.Keys) {
    $ws.Range($addr).Value = $safeUpdates[$addr]
}

foreach ($addr in $textUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textUpdates[$addr]
}

Write-Host "Applied" ($safeUpdates.Count + $textUpdates.Count) "cell updates"
